# pushing what I currently have before making big changes
#
# Adds a new "fish-dominant" worksheet (a working copy of the "2020" diet
# composition, with per_diet values reworked and a new "soybean meal "
# ingredient row appended) as the third/last sheet and leaves it as the
# active tab - mirroring the author's commit.

$wb = $excel.ActiveWorkbook

# --- Add the new sheet after the last existing one ("2020") -----------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "fish-dominant"

# --- Header row ---------------------------------------------------------
$ws.Cells.Item(1, 1).Value = "species"
$ws.Cells.Item(1, 2).Value = "production_system"
$ws.Cells.Item(1, 3).Value = "per_diet"
$ws.Cells.Item(1, 4).Value = "raw_name"

# --- Data rows ------------------------------------------------------------
# species, production_system, per_diet (blank = no value), raw_name
$rows = @(
    @("salmon", "aquaculture", $null, "soy protein concentrate"),
    @("salmon", "aquaculture", 13,    "wheat gluten"),
    @("salmon", "aquaculture", $null, "guar protein"),
    @("salmon", "aquaculture", $null, "sunflower"),
    @("salmon", "aquaculture", $null, "pea protein"),
    @("salmon", "aquaculture", 4,     "corn gluten"),
    @("salmon", "aquaculture", $null, "rapeseed oil"),
    @("salmon", "aquaculture", $null, "linseed oil"),
    @("salmon", "aquaculture", $null, "soybean oil"),
    @("salmon", "aquaculture", $null, "camelina oil"),
    @("salmon", "aquaculture", $null, "coconut oil"),
    @("salmon", "aquaculture", 7,     "wheat"),
    @("salmon", "aquaculture", 2,     "faba beans"),
    @("salmon", "aquaculture", $null, "pea flour"),
    @("salmon", "aquaculture", 32,    "fish meal, forage fish"),
    @("salmon", "aquaculture", $null, "fish meal, cut offs"),
    @("salmon", "aquaculture", 23,    "fish oil, forage fish"),
    @("salmon", "aquaculture", $null, "fish oil, cut offs"),
    @("salmon", "aquaculture", 2,     "micro ingredients"),
    @("salmon", "aquaculture", $null, "other"),
    @("salmon", "aquaculture", 17,    "soybean meal ")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    if ($null -ne $row[2]) {
        $ws.Cells.Item($r, 3).Value = $row[2]
    }
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# --- Column widths (matches the "2020" sheet it was copied from) ----------
$ws.Columns.Item(2).ColumnWidth = 16.83
$ws.Columns.Item(4).ColumnWidth = 20.5

# --- View state: new sheet becomes active with C2 selected ----------------
$ws.Activate()
$ws.Range("C2").Select()

# --- "2020" sheet: select the whole sheet (keeping old active cell) -------
$ws2020 = $wb.Worksheets.Item("2020")
$ws2020.Cells.Select()

# Re-activate the new sheet so it ends up as the selected/active tab
$ws.Activate()
